# Auto-generated edit script applying cached market-price / profit updates
# from the scheduled Hyperion pricing refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1207886.9
$ws.Range("I5").Value = 2174036.5
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 2174036.5
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -2173921.5
$ws.Range("N5").Value = -430
$ws.Range("H12").Value = 2526291
$ws.Range("I12").Value = 5682155
$ws.Range("J12").Value = 1599.8
$ws.Range("K12").Value = 5682155
$ws.Range("L12").Value = 1599.8
$ws.Range("M12").Value = -5681985
$ws.Range("N12").Value = -1939.8
$ws.Range("H39").Value = 140
$ws.Range("I39").Value = 45.25
$ws.Range("J39").Value = 224.22223
$ws.Range("K39").Value = 135.75
$ws.Range("L39").Value = 672.66669
$ws.Range("M39").Value = 160.25
$ws.Range("N39").Value = -1264.66669
$ws.Range("H64").Value = 8081.278
$ws.Range("J64").Value = 8146.706
$ws.Range("L64").Value = 8146.706
$ws.Range("N64").Value = -8642.706
$ws.Range("H67").Value = 8081.278
$ws.Range("J67").Value = 8146.706
$ws.Range("L67").Value = 8146.706
$ws.Range("N67").Value = -9862.706
$ws.Range("H103").Value = 447.14816
$ws.Range("I103").Value = 497.5909
$ws.Range("J103").Value = 225.2
$ws.Range("K103").Value = 1492.7727
$ws.Range("L103").Value = 675.5999999999999
$ws.Range("M103").Value = -906.7727
$ws.Range("N103").Value = -1847.6
$ws.Range("H115").Value = 1281.5
$ws.Range("I115").Value = 1281.5
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 3844.5
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -2277.5
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1886179.1
$ws.Range("I2").Value = 2693662.8
$ws.Range("J2").Value = 2050.5557
$ws.Range("K2").Value = 2693662.8
$ws.Range("L2").Value = 2050.5557
$ws.Range("M2").Value = -2693549.8
$ws.Range("N2").Value = -2276.5557
$ws.Range("H32").Value = 9574.581
$ws.Range("I32").Value = 5059.7393
$ws.Range("K32").Value = 5059.7393
$ws.Range("M32").Value = -4772.7393
$ws.Range("H45").Value = 6542724.5
$ws.Range("I45").Value = 13079597
$ws.Range("K45").Value = 13079597
$ws.Range("M45").Value = -13079220
$ws.Range("H74").Value = 35525.93
$ws.Range("I74").Value = 2614.762
$ws.Range("K74").Value = 2614.762
$ws.Range("M74").Value = -1740.762
$ws.Range("H77").Value = 35525.93
$ws.Range("I77").Value = 2614.762
$ws.Range("K77").Value = 13073.81
$ws.Range("M77").Value = -8705.810000000001
$ws.Range("H116").Value = 1886179.1
$ws.Range("I116").Value = 2693662.8
$ws.Range("J116").Value = 2050.5557
$ws.Range("K116").Value = 2693662.8
$ws.Range("L116").Value = 2050.5557
$ws.Range("M116").Value = -2691368.8
$ws.Range("N116").Value = -6638.5557

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1886179.1
$ws.Range("I3").Value = 2693662.8
$ws.Range("J3").Value = 2050.5557
$ws.Range("K3").Value = 2693662.8
$ws.Range("L3").Value = 2050.5557
$ws.Range("M3").Value = -2693548.8
$ws.Range("N3").Value = -2278.5557
$ws.Range("H86").Value = 2779330.8
$ws.Range("I86").Value = 3334853
$ws.Range("K86").Value = 3334853
$ws.Range("M86").Value = -3333730
$ws.Range("H89").Value = 2779330.8
$ws.Range("I89").Value = 3334853
$ws.Range("K89").Value = 16674265
$ws.Range("M89").Value = -16668649
$ws.Range("H107").Value = 2464896.5
$ws.Range("I107").Value = 3107109.5
$ws.Range("K107").Value = 3107109.5
$ws.Range("M107").Value = -3105189.5
$ws.Range("H139").Value = 123247.875
$ws.Range("I139").Value = 92500
$ws.Range("J139").Value = 133497.17
$ws.Range("K139").Value = 92500
$ws.Range("L139").Value = 133497.17
$ws.Range("M139").Value = -87360
$ws.Range("N139").Value = -143777.17

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5500
$ws.Range("I4").Value = 7000
$ws.Range("J4").Value = 4000
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 4000
$ws.Range("M4").Value = -6888
$ws.Range("N4").Value = -4224
$ws.Range("H16").Value = 1682.8572
$ws.Range("I16").Value = 961.6667
$ws.Range("J16").Value = 2223.75
$ws.Range("K16").Value = 961.6667
$ws.Range("L16").Value = 2223.75
$ws.Range("M16").Value = -674.6667
$ws.Range("N16").Value = -2797.75
$ws.Range("H31").Value = 18951.033
$ws.Range("I31").Value = 1983.8948
$ws.Range("J31").Value = 26813.854
$ws.Range("K31").Value = 1983.8948
$ws.Range("L31").Value = 26813.854
$ws.Range("M31").Value = -1688.8948
$ws.Range("N31").Value = -27403.854
$ws.Range("H34").Value = 18951.033
$ws.Range("I34").Value = 1983.8948
$ws.Range("J34").Value = 26813.854
$ws.Range("K34").Value = 1983.8948
$ws.Range("L34").Value = 26813.854
$ws.Range("M34").Value = -1781.8948
$ws.Range("N34").Value = -27217.854
$ws.Range("H113").Value = 1682.8572
$ws.Range("I113").Value = 961.6667
$ws.Range("J113").Value = 2223.75
$ws.Range("K113").Value = 961.6667
$ws.Range("L113").Value = 2223.75
$ws.Range("M113").Value = 1208.3333
$ws.Range("N113").Value = -6563.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 478714.22
$ws.Range("I2").Value = 675.8333
$ws.Range("K2").Value = 4054.9998
$ws.Range("M2").Value = -3941.9998
$ws.Range("H4").Value = 15278069
$ws.Range("I4").Value = 18661436
$ws.Range("J4").Value = 52917.5
$ws.Range("K4").Value = 55984308
$ws.Range("L4").Value = 158752.5
$ws.Range("M4").Value = -55984196
$ws.Range("N4").Value = -158976.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 527.6
$ws.Range("I107").Value = 496.57144
$ws.Range("K107").Value = 496.57144
$ws.Range("M107").Value = 1423.42856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3346666.2
$ws.Range("J2").Value = 20000
$ws.Range("L2").Value = 20000
$ws.Range("N2").Value = -20224
$ws.Range("H22").Value = 57089.25
$ws.Range("I22").Value = 64387.715
$ws.Range("K22").Value = 64387.715
$ws.Range("M22").Value = -64092.715
$ws.Range("H27").Value = 57089.25
$ws.Range("I27").Value = 64387.715
$ws.Range("K27").Value = 64387.715
$ws.Range("M27").Value = -64280.715
$ws.Range("H92").Value = 62999.75
$ws.Range("J92").Value = 62999.75
$ws.Range("L92").Value = 62999.75
$ws.Range("N92").Value = -67991.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 61969.6
$ws.Range("J46").Value = 66212
$ws.Range("L46").Value = 66212
$ws.Range("N46").Value = -66674
$ws.Range("H133").Value = 79215.336
$ws.Range("J133").Value = 79215.336
$ws.Range("L133").Value = 79215.336
$ws.Range("N133").Value = -89335.336
$ws.Range("H134").Value = 61969.6
$ws.Range("J134").Value = 66212
$ws.Range("L134").Value = 198636
$ws.Range("N134").Value = -203706
$ws.Range("H138").Value = 85229.8
$ws.Range("J138").Value = 85229.8
$ws.Range("L138").Value = 85229.8
$ws.Range("N138").Value = -95509.8
$ws.Range("H139").Value = 160993.5
$ws.Range("J139").Value = 160993.5
$ws.Range("L139").Value = 160993.5
$ws.Range("N139").Value = -171273.5
$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360
